$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Planes.pdf"
$ws.Range("B2").Value = "pdf"

$ws.Range("A3").Value = "Trains.jpg"
$ws.Range("B3").Value = "jpg"

$ws.Range("A4").Value = "Automobiles.docx"
$ws.Range("B4").Value = "docx"
